# Generate Report for Handoff
#
# A new source file "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.md" has reached
# "Ready for handoff" status (same as 633c5ece-...). It is inserted as the
# new last data row in each sheet (row 7), immediately before the
# ".localization-config" bookkeeping row which is pushed down to row 8.

$wb = $excel.ActiveWorkbook

$newFile      = "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.md"
$readyStatus  = "Ready for handoff"
$configFile   = ".localization-config"
$notLocalized = "Not to be localized"
$includeTxt   = "Include"
$ignoredTxt   = "Ignored"
$emptyDate    = "0001-01-01 00:00:00"

$zhXlf   = "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.0001c6190457a4bc7d05ec8578fa22b2ddb4258c.zh-cn.xlf"
$zhDate  = "2016-02-24 06:49:51"
$deXlf   = "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.0001c6190457a4bc7d05ec8578fa22b2ddb4258c.de-de.xlf"
$deDate  = "2016-02-24 06:50:03"

$mdCommit     = "894819da7be55e72501318f78e4fc467493c71ce"
$mdUrl        = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$newFile"
$configUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/.localization-config"
$xlfCommit    = "2c6190457a4bc7d05ec8578fa22b2ddb4258c2c"
$zhXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf"
$deXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf"

# ---------------------------------------------------------------------
# Sheet 1: "Overview" -- columns A (File Name), B (zh-cn), C (de-de)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A8").Value = $configFile
$wsOverview.Range("B8").Value = $notLocalized
$wsOverview.Range("C8").Value = $notLocalized

$wsOverview.Range("A7").Value = $newFile
$wsOverview.Range("B7").Value = $readyStatus
$wsOverview.Range("C7").Value = $readyStatus

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/894819da7be55e72501318f78e4fc467493c71ce/e2e/4ce11041-e45b-498f-8c2c-eabe1ac6898a.md", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/035bace660c33c15c88e8185bb979ad1a2bbfef7/e2e/4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.md", "", "", "4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c14d0a00ff531cf1f4ce1810cf8a6545f06b5605/e2e/7e07e4e7-cbff-4667-8c03-a4399b2129a2.md", "", "", "7e07e4e7-cbff-4667-8c03-a4399b2129a2.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c14d0a00ff531cf1f4ce1810cf8a6545f06b5605/e2e/8f271a84-b4a1-479c-bbe8-808663fd2297.md", "", "", "8f271a84-b4a1-479c-bbe8-808663fd2297.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/fa2aec7a598d2ff7c236d9cffc7dd5855780b858/e2e/633c5ece-5139-4489-b415-aa0b99d7bbec.md", "", "", "633c5ece-5139-4489-b415-aa0b99d7bbec.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), $mdUrl, "", "", $newFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A8"), $configUrl, "", "", $configFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

$wsZh.Range("A8").Value = $configFile
$wsZh.Range("B8").Value = $notLocalized
$wsZh.Range("D8").Value = $emptyDate
$wsZh.Range("G8").Value = $emptyDate
$wsZh.Range("H8").Value = $ignoredTxt

$wsZh.Range("A7").Value = $newFile
$wsZh.Range("B7").Value = $readyStatus
$wsZh.Range("C7").Value = $zhXlf
$wsZh.Range("D7").Value = $zhDate
$wsZh.Range("G7").Value = $emptyDate
$wsZh.Range("H7").Value = $includeTxt

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/894819da7be55e72501318f78e4fc467493c71ce/e2e/4ce11041-e45b-498f-8c2c-eabe1ac6898a.md", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d55a7eb561da2087d4071d4a73f4765dd8aa1c53/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4ce11041-e45b-498f-8c2c-eabe1ac6898a.0b93b2f24cd5b4391d360b2af845c418b5291816.zh-cn.xlf", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.0b93b2f24cd5b4391d360b2af845c418b5291816.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0166de96ff8005e319a317cf144e304c1813cca2/e2e/4ce11041-e45b-498f-8c2c-eabe1ac6898a.md", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b6f6b82834f45fd251086ba1b8eac7d5034e950d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4ce11041-e45b-498f-8c2c-eabe1ac6898a.0b93b2f24cd5b4391d360b2af845c418b5291816.zh-cn.xlf", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.0b93b2f24cd5b4391d360b2af845c418b5291816.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/035bace660c33c15c88e8185bb979ad1a2bbfef7/e2e/4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.md", "", "", "4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd5fdf43167231fbd7acfd63e0445ee23b59f973/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.db903e71047d7abbb01e818721e1ed7dfff3a14d.zh-cn.xlf", "", "", "4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.db903e71047d7abbb01e818721e1ed7dfff3a14d.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c14d0a00ff531cf1f4ce1810cf8a6545f06b5605/e2e/7e07e4e7-cbff-4667-8c03-a4399b2129a2.md", "", "", "7e07e4e7-cbff-4667-8c03-a4399b2129a2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca09c5deda775440cf2146696437e3325acb0f5e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7e07e4e7-cbff-4667-8c03-a4399b2129a2.4773fa57758c452be6d1778b3fe178742b112881.zh-cn.xlf", "", "", "7e07e4e7-cbff-4667-8c03-a4399b2129a2.4773fa57758c452be6d1778b3fe178742b112881.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c14d0a00ff531cf1f4ce1810cf8a6545f06b5605/e2e/8f271a84-b4a1-479c-bbe8-808663fd2297.md", "", "", "8f271a84-b4a1-479c-bbe8-808663fd2297.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca09c5deda775440cf2146696437e3325acb0f5e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8f271a84-b4a1-479c-bbe8-808663fd2297.19a759227c4efbf8e0401d2ee997133a4cb290d1.zh-cn.xlf", "", "", "8f271a84-b4a1-479c-bbe8-808663fd2297.19a759227c4efbf8e0401d2ee997133a4cb290d1.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/fa2aec7a598d2ff7c236d9cffc7dd5855780b858/e2e/633c5ece-5139-4489-b415-aa0b99d7bbec.md", "", "", "633c5ece-5139-4489-b415-aa0b99d7bbec.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1b7ca9811b5658cf993fcc527299f9e53a2e1d19/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/633c5ece-5139-4489-b415-aa0b99d7bbec.348d93e84a68523c1d12fe2ba726f3c1d928c2c4.zh-cn.xlf", "", "", "633c5ece-5139-4489-b415-aa0b99d7bbec.348d93e84a68523c1d12fe2ba726f3c1d928c2c4.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), $mdUrl, "", "", $newFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), $zhXlfUrl, "", "", $zhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A8"), $configUrl, "", "", $configFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

$wsDe.Range("A8").Value = $configFile
$wsDe.Range("B8").Value = $notLocalized
$wsDe.Range("D8").Value = $emptyDate
$wsDe.Range("G8").Value = $emptyDate
$wsDe.Range("H8").Value = $ignoredTxt

$wsDe.Range("A7").Value = $newFile
$wsDe.Range("B7").Value = $readyStatus
$wsDe.Range("C7").Value = $deXlf
$wsDe.Range("D7").Value = $deDate
$wsDe.Range("G7").Value = $emptyDate
$wsDe.Range("H7").Value = $includeTxt

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/894819da7be55e72501318f78e4fc467493c71ce/e2e/4ce11041-e45b-498f-8c2c-eabe1ac6898a.md", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a20d40ac2059ee1bdd855f7b2a3075731d6b1fc0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4ce11041-e45b-498f-8c2c-eabe1ac6898a.0b93b2f24cd5b4391d360b2af845c418b5291816.de-de.xlf", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.0b93b2f24cd5b4391d360b2af845c418b5291816.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b7824b29487544aeb3ed800cef6f2c9608a3ce34/e2e/4ce11041-e45b-498f-8c2c-eabe1ac6898a.md", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6a4e632325a71625a02e1c87baddbc7f7723908a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4ce11041-e45b-498f-8c2c-eabe1ac6898a.0b93b2f24cd5b4391d360b2af845c418b5291816.de-de.xlf", "", "", "4ce11041-e45b-498f-8c2c-eabe1ac6898a.0b93b2f24cd5b4391d360b2af845c418b5291816.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/035bace660c33c15c88e8185bb979ad1a2bbfef7/e2e/4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.md", "", "", "4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8fceaef1fb7a3d6e787f99fb7465bb302d478173/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.db903e71047d7abbb01e818721e1ed7dfff3a14d.de-de.xlf", "", "", "4fb2799b-b9ff-4ed0-8e7e-faee65df21cf.db903e71047d7abbb01e818721e1ed7dfff3a14d.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c14d0a00ff531cf1f4ce1810cf8a6545f06b5605/e2e/7e07e4e7-cbff-4667-8c03-a4399b2129a2.md", "", "", "7e07e4e7-cbff-4667-8c03-a4399b2129a2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3823035973a26a36ef0dfd79cefb8abafafcf9cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7e07e4e7-cbff-4667-8c03-a4399b2129a2.4773fa57758c452be6d1778b3fe178742b112881.de-de.xlf", "", "", "7e07e4e7-cbff-4667-8c03-a4399b2129a2.4773fa57758c452be6d1778b3fe178742b112881.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c14d0a00ff531cf1f4ce1810cf8a6545f06b5605/e2e/8f271a84-b4a1-479c-bbe8-808663fd2297.md", "", "", "8f271a84-b4a1-479c-bbe8-808663fd2297.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3823035973a26a36ef0dfd79cefb8abafafcf9cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8f271a84-b4a1-479c-bbe8-808663fd2297.19a759227c4efbf8e0401d2ee997133a4cb290d1.de-de.xlf", "", "", "8f271a84-b4a1-479c-bbe8-808663fd2297.19a759227c4efbf8e0401d2ee997133a4cb290d1.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/fa2aec7a598d2ff7c236d9cffc7dd5855780b858/e2e/633c5ece-5139-4489-b415-aa0b99d7bbec.md", "", "", "633c5ece-5139-4489-b415-aa0b99d7bbec.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c86341d40f69f92d51b943680e039c8b2cb0fd09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/633c5ece-5139-4489-b415-aa0b99d7bbec.348d93e84a68523c1d12fe2ba726f3c1d928c2c4.de-de.xlf", "", "", "633c5ece-5139-4489-b415-aa0b99d7bbec.348d93e84a68523c1d12fe2ba726f3c1d928c2c4.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), $mdUrl, "", "", $newFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), $deXlfUrl, "", "", $deXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A8"), $configUrl, "", "", $configFile) | Out-Null
